# Trade #55 closed at 2026-02-17 21:10:14 - unknown UNKNOWN +0.000%
# Also a new OPEN trade (#116) was logged at 21:10:07.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.92
$summary.Range("B4").Value = 0.72
$summary.Range("B5").Value = 0.17
$summary.Range("B6").Value = 83
$summary.Range("B8").Value = 32
$summary.Range("B9").Value = 46.99

# ---------------------------------------------------------------------
# Sheet: Strategy Status (MarketMaking row = row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.92
$status.Range("D5").Value = 50
$status.Range("E5").Value = 0.61
$status.Range("F5").Value = 0.92
$status.Range("G5").Value = 50

# ---------------------------------------------------------------------
# Sheet: All Trades
#   - row 84 (Trade #83) : closed with an early exit
#   - row 117 (Trade #116): new OPEN trade appended
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Range("G84").Value = 0.3
$allTrades.Range("H84").Value = "CLOSED"
$allTrades.Range("I84").Value = -30.2326
$allTrades.Range("J84").Value = -0.13
$allTrades.Range("K84").Value = 100.92
$allTrades.Range("L84").Value = "early_exit"
$allTrades.Range("M84").Value = 0.13

$allTrades.Range("A117").Value = 116
$allTrades.Range("B117").NumberFormat = "@"
$allTrades.Range("B117").Value = "2026-02-17"
$allTrades.Range("C117").Value = "21:10:07"
$allTrades.Range("D117").Value = "MarketMaking"
$allTrades.Range("E117").Value = "DOWN"
$allTrades.Range("F117").Value = 0.43
$allTrades.Range("G117").Value = ""
$allTrades.Range("H117").Value = "OPEN"
$allTrades.Range("I117").Value = 0
$allTrades.Range("J117").Value = 0
$allTrades.Range("K117").Value = 101.0514872031006
$allTrades.Range("L117").Value = ""
$allTrades.Range("M117").Value = 0
$allTrades.Range("N117").Value = 0
$allTrades.Range("O117").Value = 0
$allTrades.Range("P117").Value = 0.6
$allTrades.Range("Q117").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# Sheet: MarketMaking
#   - row 51 (Trade #83) : closed with an early exit
#   - row 84 (Trade #116): new OPEN trade appended
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

$mm.Range("G51").Value = 0.3
$mm.Range("H51").Value = "CLOSED"
$mm.Range("I51").Value = -30.2326
$mm.Range("J51").Value = -0.13
$mm.Range("K51").Value = 100.92
$mm.Range("P51").Value = "early_exit"
$mm.Range("Q51").Value = 0.13

$mm.Range("A84").Value = 116
$mm.Range("B84").NumberFormat = "@"
$mm.Range("B84").Value = "2026-02-17"
$mm.Range("C84").Value = "21:10:07"
$mm.Range("D84").Value = "MarketMaking"
$mm.Range("E84").Value = "DOWN"
$mm.Range("F84").Value = 0.43
$mm.Range("G84").Value = ""
$mm.Range("H84").Value = "OPEN"
$mm.Range("I84").Value = 0
$mm.Range("J84").Value = 0
$mm.Range("K84").Value = 101.0514872031006
$mm.Range("L84").Value = 0
$mm.Range("M84").Value = 0
$mm.Range("N84").Value = 0.6
$mm.Range("O84").Value = "Normal spread capture: 19600 bps"
$mm.Range("P84").Value = ""
$mm.Range("Q84").Value = 0
